# Apply the "Trade #66 closed" update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet - update aggregate stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.15   # Current Capital
$summary.Range("B4").Value = 0.14      # Total P&L $
$summary.Range("B5").Value = 0.04      # Total P&L %
$summary.Range("B6").Value = 66        # Total Trades
$summary.Range("B8").Value = 22        # Losing Trades
$summary.Range("B9").Value = 43.94     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - update MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.15     # Capital
$status.Range("D4").Value = 66         # Trades
$status.Range("E4").Value = 0.14       # P&L $
$status.Range("F4").Value = 0.15       # P&L %
$status.Range("G4").Value = 43.94      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append new closed trade (#66) as row 67 on "All Trades" and
#    "MarketMaking" sheets (both keep an identical trade log).
# ---------------------------------------------------------------------
function Add-TradeRow67($ws) {
    $ws.Cells.Item(67, 1).Value = 66

    # Date/Time columns hold plain text in this workbook (not real Excel
    # dates), so force a text format before writing then clear the
    # leftover number format so no stray style is left behind.
    $ws.Cells.Item(67, 2).NumberFormat = "@"
    $ws.Cells.Item(67, 2).Value = "2026-02-17"
    $ws.Cells.Item(67, 2).Style = "Normal"

    $ws.Cells.Item(67, 3).NumberFormat = "@"
    $ws.Cells.Item(67, 3).Value = "12:53:38"
    $ws.Cells.Item(67, 3).Style = "Normal"

    $ws.Cells.Item(67, 4).Value = "MarketMaking"
    $ws.Cells.Item(67, 5).Value = "DOWN"
    $ws.Cells.Item(67, 6).Value = 0.92
    $ws.Cells.Item(67, 7).Value = 0.91
    $ws.Cells.Item(67, 8).Value = "CLOSED"
    $ws.Cells.Item(67, 9).Value = -1.087
    $ws.Cells.Item(67, 10).Value = -0.01
    $ws.Cells.Item(67, 11).Value = 100.15
    $ws.Cells.Item(67, 12).Value = 0
    $ws.Cells.Item(67, 13).Value = 0
    $ws.Cells.Item(67, 14).Value = 0.6
    $ws.Cells.Item(67, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(67, 16).Value = "early_exit"
    $ws.Cells.Item(67, 17).Value = 0.14
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow67 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow67 $marketMaking
